$d = $word.ActiveDocument

# --- Change 1: Replace "[TITLE]" run with bold title text, and " investigates" with " studies" ---

# Find the "[TITLE]" text and replace it with the bold title
$range1 = $d.Content
$range1.Find.Execute("[TITLE]", $false, $false, $false, $false, $false, $true, 1, $false, "Investigating the Impact of Augmented Reality and BIM on Retrofitting Training for Non-experts", 2) | Out-Null

# Now make that replaced text bold. Find it again to get its range.
$range1b = $d.Content
$range1b.Find.Execute("Investigating the Impact of Augmented Reality and BIM on Retrofitting Training for Non-experts", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$range1b.Bold = 1

# Replace " investigates" with " studies"
$range2 = $d.Content
$range2.Find.Execute(" investigates", $true, $false, $false, $false, $false, $true, 1, $false, " studies", 2) | Out-Null

# --- Change 2: Insert "s" after "effort" (before " required to complete the installation") ---

$range3 = $d.Content
$range3.Find.Execute("effort required to complete the installation", $true, $false, $false, $false, $false, $true, 1, $false, "efforts required to complete the installation", 2) | Out-Null
